$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the empty E92/F92 cells (row 92 should end at column D)
$ws.Range("E92:F92").ClearContents()

# Ensure columns A:D are treated as Text so values like dates/numbers are not auto-converted
$ws.Range("A93:D155").NumberFormat = "@"

# Populate new rows 93-155 with data (columns A-D)
$ws.Range("A93").Value = "1"
$ws.Range("B93").Value = "a"
$ws.Range("C93").Value = "2024-10-05"
$ws.Range("D93").Value = "P"

$ws.Range("A94").Value = "2"
$ws.Range("B94").Value = "b"
$ws.Range("C94").Value = "2024-10-05"
$ws.Range("D94").Value = "P"

$ws.Range("A95").Value = "3"
$ws.Range("B95").Value = "c"
$ws.Range("C95").Value = "2024-10-05"
$ws.Range("D95").Value = "P"

$ws.Range("A96").Value = "123123"
$ws.Range("B96").Value = "aluno 1000"
$ws.Range("C96").Value = "2024-10-05"
$ws.Range("D96").Value = "P"

$ws.Range("A97").Value = "7676887"
$ws.Range("B97").Value = "aluno 2000"
$ws.Range("C97").Value = "2024-10-05"
$ws.Range("D97").Value = "P"

$ws.Range("A98").Value = "1"
$ws.Range("B98").Value = "a"
$ws.Range("C98").Value = "2024-10-05"
$ws.Range("D98").Value = "P"

$ws.Range("A99").Value = "2"
$ws.Range("B99").Value = "b"
$ws.Range("C99").Value = "2024-10-05"
$ws.Range("D99").Value = "P"

$ws.Range("A100").Value = "3"
$ws.Range("B100").Value = "c"
$ws.Range("C100").Value = "2024-10-05"
$ws.Range("D100").Value = "P"

$ws.Range("A101").Value = "123123"
$ws.Range("B101").Value = "aluno 1000"
$ws.Range("C101").Value = "2024-10-05"
$ws.Range("D101").Value = "P"

$ws.Range("A102").Value = "7676887"
$ws.Range("B102").Value = "aluno 2000"
$ws.Range("C102").Value = "2024-10-05"
$ws.Range("D102").Value = "P"

$ws.Range("A103").Value = "1"
$ws.Range("B103").Value = "a"
$ws.Range("C103").Value = "2024-01-01"
$ws.Range("D103").Value = "A"

$ws.Range("A104").Value = "2"
$ws.Range("B104").Value = "b"
$ws.Range("C104").Value = "2024-01-01"
$ws.Range("D104").Value = "A"

$ws.Range("A105").Value = "3"
$ws.Range("B105").Value = "c"
$ws.Range("C105").Value = "2024-01-01"
$ws.Range("D105").Value = "A"

$ws.Range("A106").Value = "123123"
$ws.Range("B106").Value = "aluno 1000"
$ws.Range("C106").Value = "2024-01-01"
$ws.Range("D106").Value = "A"

$ws.Range("A107").Value = "7676887"
$ws.Range("B107").Value = "aluno 2000"
$ws.Range("C107").Value = "2024-01-01"
$ws.Range("D107").Value = "A"

$ws.Range("A108").Value = "1"
$ws.Range("B108").Value = "a"
$ws.Range("C108").Value = "2024-01-02"
$ws.Range("D108").Value = "P"

$ws.Range("A109").Value = "2"
$ws.Range("B109").Value = "b"
$ws.Range("C109").Value = "2024-01-02"
$ws.Range("D109").Value = "P"

$ws.Range("A110").Value = "3"
$ws.Range("B110").Value = "c"
$ws.Range("C110").Value = "2024-01-02"
$ws.Range("D110").Value = "P"

$ws.Range("A111").Value = "123123"
$ws.Range("B111").Value = "aluno 1000"
$ws.Range("C111").Value = "2024-01-02"
$ws.Range("D111").Value = "P"

$ws.Range("A112").Value = "7676887"
$ws.Range("B112").Value = "aluno 2000"
$ws.Range("C112").Value = "2024-01-02"
$ws.Range("D112").Value = "P"

$ws.Range("A113").Value = "1"
$ws.Range("B113").Value = "a"
$ws.Range("C113").Value = "2024-01-03"
$ws.Range("D113").Value = "P"

$ws.Range("A114").Value = "2"
$ws.Range("B114").Value = "b"
$ws.Range("C114").Value = "2024-01-03"
$ws.Range("D114").Value = "P"

$ws.Range("A115").Value = "3"
$ws.Range("B115").Value = "c"
$ws.Range("C115").Value = "2024-01-03"
$ws.Range("D115").Value = "P"

$ws.Range("A116").Value = "123123"
$ws.Range("B116").Value = "aluno 1000"
$ws.Range("C116").Value = "2024-01-03"
$ws.Range("D116").Value = "A"

$ws.Range("A117").Value = "7676887"
$ws.Range("B117").Value = "aluno 2000"
$ws.Range("C117").Value = "2024-01-03"
$ws.Range("D117").Value = "A"

$ws.Range("A118").Value = "1"
$ws.Range("B118").Value = "a"
$ws.Range("C118").Value = "2024-02-01"
$ws.Range("D118").Value = "P"

$ws.Range("A119").Value = "2"
$ws.Range("B119").Value = "b"
$ws.Range("C119").Value = "2024-02-01"
$ws.Range("D119").Value = "P"

$ws.Range("A120").Value = "3"
$ws.Range("B120").Value = "c"
$ws.Range("C120").Value = "2024-02-01"
$ws.Range("D120").Value = "P"

$ws.Range("A121").Value = "123123"
$ws.Range("B121").Value = "aluno 1000"
$ws.Range("C121").Value = "2024-02-01"
$ws.Range("D121").Value = "P"

$ws.Range("A122").Value = "7676887"
$ws.Range("B122").Value = "aluno 2000"
$ws.Range("C122").Value = "2024-02-01"
$ws.Range("D122").Value = "A"

$ws.Range("A123").Value = "1"
$ws.Range("B123").Value = "a"
$ws.Range("C123").Value = "2024-02-02"
$ws.Range("D123").Value = "A"

$ws.Range("A124").Value = "2"
$ws.Range("B124").Value = "b"
$ws.Range("C124").Value = "2024-02-02"
$ws.Range("D124").Value = "A"

$ws.Range("A125").Value = "3"
$ws.Range("B125").Value = "c"
$ws.Range("C125").Value = "2024-02-02"
$ws.Range("D125").Value = "A"

$ws.Range("A126").Value = "123123"
$ws.Range("B126").Value = "aluno 1000"
$ws.Range("C126").Value = "2024-02-02"
$ws.Range("D126").Value = "A"

$ws.Range("A127").Value = "7676887"
$ws.Range("B127").Value = "aluno 2000"
$ws.Range("C127").Value = "2024-02-02"
$ws.Range("D127").Value = "A"

$ws.Range("A128").Value = "1"
$ws.Range("B128").Value = "a"
$ws.Range("C128").Value = "2024-03-01"
$ws.Range("D128").Value = "P"

$ws.Range("A129").Value = "2"
$ws.Range("B129").Value = "b"
$ws.Range("C129").Value = "2024-03-01"
$ws.Range("D129").Value = "P"

$ws.Range("A130").Value = "3"
$ws.Range("B130").Value = "c"
$ws.Range("C130").Value = "2024-03-01"
$ws.Range("D130").Value = "P"

$ws.Range("A131").Value = "123123"
$ws.Range("B131").Value = "aluno 1000"
$ws.Range("C131").Value = "2024-03-01"
$ws.Range("D131").Value = "P"

$ws.Range("A132").Value = "7676887"
$ws.Range("B132").Value = "aluno 2000"
$ws.Range("C132").Value = "2024-03-01"
$ws.Range("D132").Value = "P"

$ws.Range("A133").Value = "1"
$ws.Range("B133").Value = "a"
$ws.Range("C133").Value = "2024-03-02"
$ws.Range("D133").Value = "A"

$ws.Range("A134").Value = "2"
$ws.Range("B134").Value = "b"
$ws.Range("C134").Value = "2024-03-02"
$ws.Range("D134").Value = "A"

$ws.Range("A135").Value = "3"
$ws.Range("B135").Value = "c"
$ws.Range("C135").Value = "2024-03-02"
$ws.Range("D135").Value = "A"

$ws.Range("A136").Value = "123123"
$ws.Range("B136").Value = "aluno 1000"
$ws.Range("C136").Value = "2024-03-02"
$ws.Range("D136").Value = "A"

$ws.Range("A137").Value = "7676887"
$ws.Range("B137").Value = "aluno 2000"
$ws.Range("C137").Value = "2024-03-02"
$ws.Range("D137").Value = "A"

$ws.Range("A138").Value = "1"
$ws.Range("B138").Value = "a"
$ws.Range("C138").Value = "2024-03-04"
$ws.Range("D138").Value = "P"

$ws.Range("A139").Value = "2"
$ws.Range("B139").Value = "b"
$ws.Range("C139").Value = "2024-03-04"
$ws.Range("D139").Value = "P"

$ws.Range("A140").Value = "3"
$ws.Range("B140").Value = "c"
$ws.Range("C140").Value = "2024-03-04"
$ws.Range("D140").Value = "P"

$ws.Range("A141").Value = "123123"
$ws.Range("B141").Value = "aluno 1000"
$ws.Range("C141").Value = "2024-03-04"
$ws.Range("D141").Value = "P"

$ws.Range("A142").Value = "7676887"
$ws.Range("B142").Value = "aluno 2000"
$ws.Range("C142").Value = "2024-03-04"
$ws.Range("D142").Value = "P"

$ws.Range("A143").Value = "1"
$ws.Range("B143").Value = "a"
$ws.Range("C143").Value = "2024-03-05"
$ws.Range("D143").Value = "P"

$ws.Range("A144").Value = "2"
$ws.Range("B144").Value = "b"
$ws.Range("C144").Value = "2024-03-05"
$ws.Range("D144").Value = "P"

$ws.Range("A145").Value = "3"
$ws.Range("B145").Value = "c"
$ws.Range("C145").Value = "2024-03-05"
$ws.Range("D145").Value = "P"

$ws.Range("A146").Value = "123123"
$ws.Range("B146").Value = "aluno 1000"
$ws.Range("C146").Value = "2024-03-05"
$ws.Range("D146").Value = "P"

$ws.Range("A147").Value = "7676887"
$ws.Range("B147").Value = "aluno 2000"
$ws.Range("C147").Value = "2024-03-05"
$ws.Range("D147").Value = "P"

$ws.Range("A148").Value = "1"
$ws.Range("B148").Value = "a"
$ws.Range("C148").Value = "2024-03-05"
$ws.Range("D148").Value = "A"

$ws.Range("A149").Value = "2"
$ws.Range("B149").Value = "b"
$ws.Range("C149").Value = "2024-03-05"
$ws.Range("D149").Value = "A"

$ws.Range("A150").Value = "3"
$ws.Range("B150").Value = "c"
$ws.Range("C150").Value = "2024-03-05"
$ws.Range("D150").Value = "A"

$ws.Range("A151").Value = "123123"
$ws.Range("B151").Value = "aluno 1000"
$ws.Range("C151").Value = "2024-03-05"
$ws.Range("D151").Value = "A"

$ws.Range("A152").Value = "7676887"
$ws.Range("B152").Value = "aluno 2000"
$ws.Range("C152").Value = "2024-03-05"
$ws.Range("D152").Value = "A"

$ws.Range("A153").Value = "1"
$ws.Range("B153").Value = "Pai do leonardo"
$ws.Range("C153").Value = "2024-01-01"
$ws.Range("D153").Value = "P"

$ws.Range("A154").Value = "1"
$ws.Range("B154").Value = "Pai do leonardo"
$ws.Range("C154").Value = "2024-01-02"
$ws.Range("D154").Value = "P"

$ws.Range("A155").Value = "1"
$ws.Range("B155").Value = "Pai do leonardo"
$ws.Range("C155").Value = "2024-01-03"
$ws.Range("D155").Value = "A"

# Row 155 also has empty E/F cells, matching the old pattern from row 92.
# Materialize them first with a placeholder, then clear to empty string so the cells persist as empty.
$ws.Range("E155:F155").NumberFormat = "@"
$ws.Range("E155").Value = "x"
$ws.Range("F155").Value = "x"
$ws.Range("E155").Value = ""
$ws.Range("F155").Value = ""
